$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.059.65"
$ws.Range("E2").Value = "  -0.15%  "

$ws.Range("D3").Value = "1.620.85"
$ws.Range("E3").Value = "  -1.06%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.87"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.34%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.514"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.64%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("E8").Value = "  +0.31%  "

$ws.Range("E9").Value = "  -1.45%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0842"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.54%  "

$ws.Range("D12").Value = "1.848.73"
$ws.Range("E12").Value = "  -1.01%  "

$ws.Range("D13").Value = "1.616.64"
$ws.Range("E13").Value = "  -1.30%  "

$ws.Range("E14").Value = "  -0.35%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.537"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.75%  "

$ws.Range("D16").Value = "27.030.03"
$ws.Range("E16").Value = "  -0.30%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.41%  "

$ws.Range("D18").Value = "0.0₃0737"
$ws.Range("E18").Value = "  -0.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "213.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.58%  "

$ws.Range("E20").Value = "  -0.11%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.83"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.97%  "

$ws.Range("E22").Value = "  -1.93%  "

$ws.Range("E23").Value = "  -7.79%  "

$ws.Range("E24").Value = "  -1.22%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.91%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.86%  "

$ws.Range("E28").Value = "  -3.38%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.46"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.45%  "

$ws.Range("E30").Value = "  +0.33%  "

$ws.Range("E31").Value = "  -1.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.08%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.688"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +26.98%  "

$ws.Range("E34").Value = "  -0.65%  "

$ws.Range("D35").Value = "1.337.87"
$ws.Range("E35").Value = "  +2.27%  "

$ws.Range("E36").Value = "  -1.21%  "

$ws.Range("E37").Value = "  -0.69%  "

$ws.Range("E38").Value = "  -0.20%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.839"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.05%  "

$ws.Range("E40").Value = "  -0.08%  "

$ws.Range("E41").Value = "  +0.32%  "

$ws.Range("E42").Value = "  -1.84%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.12%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "63.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.03%  "

$ws.Range("D45").Value = "1.759.39"
$ws.Range("E45").Value = "  -1.04%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "89.78"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.62%  "

$ws.Range("E47").Value = "  +2.36%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.848"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +26.87%  "

$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0999"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.19%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0512"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.03%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.43%  "
